# LV Contacts - 11 July 2024
# Add a "ContactName" column to the Contact sheet and make Contact the
# active/selected sheet (previously AddRelationship was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# Make "Contact" the active sheet/tab (moves tabSelected from AddRelationship
# and updates workbook.xml's bookViews/activeTab automatically).
$ws.Activate()

# New column F: header + single data row.
$ws.Range("F1").Value = "ContactName"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = "Test LVContact"

# Match the column width recorded for the new column.
$ws.Columns.Item(6).ColumnWidth = 13

# Move the active selection to F8, as in the edited workbook.
$ws.Range("F8").Select()
